$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column Q, row 4: header year 2020 - same style as P4 (s=15)
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

# Column Q, row 5: value 53.2 - same style as P5 (s=16)
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 53.2

# Column Q, row 6: value 23.2 - same style as P6 (s=31)
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 23.2

# Column Q, row 7: value 10 - like P7 but with a one-decimal number format
# (new style xf, since the value needs to render as "10.0")
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").NumberFormat = "0.0"
$ws.Range("Q7").Value = 10

# Column Q, row 8: value 20 - like P8 but with a one-decimal number format
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").NumberFormat = "0.0"
$ws.Range("Q8").Value = 20

# Match the author's final selection state
$ws.Range("P9").Select()
